$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Index_of_functions")

# Fix typo in the title cell
$ws.Range("A1").Value = "Index of Functions"

# Insert a new row above row 279 (shifts existing rows down)
$ws.Rows.Item(279).Insert()

# Fill in the new row's data: function name "readLines" used in Week 8b (column L)
$ws.Range("A279").Value = "readLines"
$ws.Range("L279").Value = "X"

# Update selection / view to match final state
$ws.Range("A279").Select()
